$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -1.277764249695792
$ws.Range("D2").Value = 0.2146515468495707

$ws.Range("C3").Value = -1.090634122122024
$ws.Range("D3").Value = 0.2872355508068913

$ws.Range("C4").Value = -1.637746939638575
$ws.Range("D4").Value = 0.1156996624296283

$ws.Range("C5").Value = -0.2135469831052462
$ws.Range("D5").Value = 0.8328669779615083

$ws.Range("C6").Value = 0.4311312426375424
$ws.Range("D6").Value = 0.6705679721484286

$ws.Range("C7").Value = 0.06793892220197796
$ws.Range("D7").Value = 0.9464479503244618

$ws.Range("C8").Value = 0.9176137987026109
$ws.Range("D8").Value = 0.3687718901240213

$ws.Range("C9").Value = -0.5064365978694232
$ws.Range("D9").Value = 0.6175880866523229

$ws.Range("C10").Value = 0.7678900580027761
$ws.Range("D10").Value = 0.4507165179005195

$ws.Range("C11").Value = 1.388275669633499
$ws.Range("D11").Value = 0.1789517303192221
